$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.801.36'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '1.613.42'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.01'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("D12").Value = '1.840.35'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.621.29'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.533'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").Value = '26.825.54'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.87'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.43'
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("E23").Value = '  -6.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.42'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.48'
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.38'
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.707'
$ws.Range("E33").Value = '  +30.47%  '
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").Value = '1.322.59'
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.53'
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.794'
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("E42").Value = '  -2.14%  '
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.44'
$ws.Range("E44").Value = '  +2.43%  '
$ws.Range("D45").Value = '1.752.33'
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.18'
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.815'
$ws.Range("E48").Value = '  +3.72%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0979'
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("E51").Value = '  -0.60%  '
